$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 12.58410635944474
$ws.Cells.Item(2, 3).Value = 7.246377296395987
$ws.Cells.Item(2, 4).Value = 7.804520551221182
$ws.Cells.Item(2, 5).Value = 13.01364591543346
$ws.Cells.Item(2, 6).Value = 38.81469517463432
$ws.Cells.Item(2, 9).Value = 29.36583346004903
$ws.Cells.Item(2, 10).Value = 10.3624905419177
$ws.Cells.Item(2, 11).Value = 10.06824227607598
$ws.Cells.Item(2, 12).Value = 11.13782472131573
$ws.Cells.Item(2, 13).Value = 14.98491546710326
$ws.Cells.Item(2, 15).Value = 30.07246463931393
$ws.Cells.Item(3, 2).Value = 12.38499528303869
$ws.Cells.Item(3, 3).Value = 7.194066479136922
$ws.Cells.Item(3, 4).Value = 7.793294854944984
$ws.Cells.Item(3, 5).Value = 13.03679928334496
$ws.Cells.Item(3, 6).Value = 38.92015916069355
$ws.Cells.Item(3, 9).Value = 29.46717665757958
$ws.Cells.Item(3, 10).Value = 10.37960624037702
$ws.Cells.Item(3, 11).Value = 9.922934195093026
$ws.Cells.Item(3, 12).Value = 11.14571950887441
$ws.Cells.Item(3, 13).Value = 14.95777706981462
$ws.Cells.Item(3, 15).Value = 30.16798427002515
$ws.Cells.Item(4, 2).Value = 12.26296533776107
$ws.Cells.Item(4, 3).Value = 7.161339742216141
$ws.Cells.Item(4, 4).Value = 7.78732804598067
$ws.Cells.Item(4, 5).Value = 13.05225367842457
$ws.Cells.Item(4, 6).Value = 38.99187148733174
$ws.Cells.Item(4, 9).Value = 29.53377498106239
$ws.Cells.Item(4, 10).Value = 10.39067825499215
$ws.Cells.Item(4, 11).Value = 9.83396397963582
$ws.Cells.Item(4, 12).Value = 11.15170351253015
$ws.Cells.Item(4, 13).Value = 14.94280531527293
$ws.Cells.Item(4, 15).Value = 30.2315512850086
$ws.Cells.Item(5, 2).Value = 12.21335399516979
$ws.Cells.Item(5, 3).Value = 7.147855362867006
$ws.Cells.Item(5, 4).Value = 7.78513124376901
$ws.Cells.Item(5, 5).Value = 13.05886336288563
$ws.Cells.Item(5, 6).Value = 39.02284296610618
$ws.Cells.Item(5, 9).Value = 29.56201461733681
$ws.Cells.Item(5, 10).Value = 10.3953321401477
$ws.Cells.Item(5, 11).Value = 9.797811637268573
$ws.Cells.Item(5, 12).Value = 11.15442854263811
$ws.Cells.Item(5, 13).Value = 14.93713375477335
$ws.Cells.Item(5, 15).Value = 30.25869163816906
$ws.Cells.Item(6, 2).Value = 12.20512489676317
$ws.Cells.Item(6, 3).Value = 7.14560749685018
$ws.Cells.Item(6, 4).Value = 7.784780696690818
$ws.Cells.Item(6, 5).Value = 13.05997974971458
$ws.Cells.Item(6, 6).Value = 39.02809129929351
$ws.Cells.Item(6, 9).Value = 29.5667702597006
$ws.Cells.Item(6, 10).Value = 10.39611349945204
$ws.Cells.Item(6, 11).Value = 9.791816034705899
$ws.Cells.Item(6, 12).Value = 11.15489835550075
$ws.Cells.Item(6, 13).Value = 14.93621805980376
$ws.Cells.Item(6, 15).Value = 30.26327293277636
$ws.Cells.Item(7, 2).Value = 12.26229570876998
$ws.Cells.Item(7, 3).Value = 7.161158479869568
$ws.Cells.Item(7, 4).Value = 7.787297466309898
$ws.Cells.Item(7, 5).Value = 13.0523415553601
$ws.Cells.Item(7, 6).Value = 38.99228210321926
$ws.Cells.Item(7, 9).Value = 29.53415137487462
$ws.Cells.Item(7, 10).Value = 10.39074044363255
$ws.Cells.Item(7, 11).Value = 9.833475942193987
$ws.Cells.Item(7, 12).Value = 11.15173910230885
$ws.Cells.Item(7, 13).Value = 14.94272708194982
$ws.Cells.Item(7, 15).Value = 30.2319123040502
$ws.Cells.Item(8, 2).Value = 12.51544212678321
$ws.Cells.Item(8, 3).Value = 7.22846954097021
$ws.Cells.Item(8, 4).Value = 7.800458962107403
$ws.Cells.Item(8, 5).Value = 13.02137255627061
$ws.Cells.Item(8, 6).Value = 38.84961449975892
$ws.Cells.Item(8, 9).Value = 29.3998693023023
$ws.Cells.Item(8, 10).Value = 10.36827544588713
$ws.Cells.Item(8, 11).Value = 10.01811279151428
$ws.Cells.Item(8, 12).Value = 11.14031137705507
$ws.Cells.Item(8, 13).Value = 14.97520973029116
$ws.Cells.Item(8, 15).Value = 30.10437900006792
$ws.Cells.Item(9, 2).Value = 13.01103611571344
$ws.Cells.Item(9, 3).Value = 7.355435849690571
$ws.Cells.Item(9, 4).Value = 7.83352618276894
$ws.Cells.Item(9, 5).Value = 12.97044188646523
$ws.Cells.Item(9, 6).Value = 38.62509793578244
$ws.Cells.Item(9, 9).Value = 29.17121445272219
$ws.Cells.Item(9, 10).Value = 10.3286696768517
$ws.Cells.Item(9, 11).Value = 10.38037951472757
$ws.Cells.Item(9, 12).Value = 11.12688851664167
$ws.Cells.Item(9, 13).Value = 15.05213288979572
$ws.Cells.Item(9, 15).Value = 29.89331760710143
$ws.Cells.Item(10, 2).Value = 13.37122926788395
$ws.Cells.Item(10, 3).Value = 7.445377466505978
$ws.Cells.Item(10, 4).Value = 7.862120667689434
$ws.Cells.Item(10, 5).Value = 12.93896369370754
$ws.Cells.Item(10, 6).Value = 38.49389752951046
$ws.Cells.Item(10, 9).Value = 29.02432052228539
$ws.Cells.Item(10, 10).Value = 10.3022571440571
$ws.Cells.Item(10, 11).Value = 10.64430490157427
$ws.Cells.Item(10, 12).Value = 11.12246242670251
$ws.Cells.Item(10, 13).Value = 15.1164368620734
$ws.Cells.Item(10, 15).Value = 29.76206111843607
$ws.Cells.Item(11, 2).Value = 13.53354713994566
$ws.Cells.Item(11, 3).Value = 7.485511803575077
$ws.Cells.Item(11, 4).Value = 7.876033870330066
$ws.Cells.Item(11, 5).Value = 12.92592621830805
$ws.Cells.Item(11, 6).Value = 38.44155171748027
$ws.Cells.Item(11, 9).Value = 28.96206842909629
$ws.Cells.Item(11, 10).Value = 10.29081911561174
$ws.Cells.Item(11, 11).Value = 10.76340594084596
$ws.Cells.Item(11, 12).Value = 11.12161933640091
$ws.Cells.Item(11, 13).Value = 15.14731861439834
$ws.Cells.Item(11, 15).Value = 29.70752257859951
$ws.Cells.Item(12, 2).Value = 13.5947368301246
$ws.Cells.Item(12, 3).Value = 7.500592591093572
$ws.Cells.Item(12, 4).Value = 7.881429985092188
$ws.Cells.Item(12, 5).Value = 12.92117305059115
$ws.Cells.Item(12, 6).Value = 38.42278562840909
$ws.Cells.Item(12, 9).Value = 28.9391519838523
$ws.Cells.Item(12, 10).Value = 10.28657041556406
$ws.Cells.Item(12, 11).Value = 10.8083297985658
$ws.Cells.Item(12, 12).Value = 11.12146746386483
$ws.Cells.Item(12, 13).Value = 15.15924131839292
$ws.Cells.Item(12, 15).Value = 29.68761411169611
$ws.Cells.Item(13, 2).Value = 13.5815716973872
$ws.Cells.Item(13, 3).Value = 7.497349965159038
$ws.Cells.Item(13, 4).Value = 7.880262209591652
$ws.Cells.Item(13, 5).Value = 12.92218856232796
$ws.Cells.Item(13, 6).Value = 38.42678025633914
$ws.Cells.Item(13, 9).Value = 28.94405822970468
$ws.Cells.Item(13, 10).Value = 10.28748177974767
$ws.Cells.Item(13, 11).Value = 10.79866312633434
$ws.Cells.Item(13, 12).Value = 11.12149274265138
$ws.Cells.Item(13, 13).Value = 15.15666348597997
$ws.Cells.Item(13, 15).Value = 29.69186865386902
$ws.Cells.Item(14, 2).Value = 13.53858710596842
$ws.Cells.Item(14, 3).Value = 7.486754886399248
$ws.Cells.Item(14, 4).Value = 7.87647527284399
$ws.Cells.Item(14, 5).Value = 12.92553149091051
$ws.Cells.Item(14, 6).Value = 38.43998664871506
$ws.Cells.Item(14, 9).Value = 28.96016991166631
$ws.Cells.Item(14, 10).Value = 10.29046791847654
$ws.Cells.Item(14, 11).Value = 10.76710562213044
$ws.Cells.Item(14, 12).Value = 11.12160349241713
$ws.Cells.Item(14, 13).Value = 15.14829495826426
$ws.Cells.Item(14, 15).Value = 29.70586978119072
$ws.Cells.Item(15, 2).Value = 13.51222017285056
$ws.Cells.Item(15, 3).Value = 7.480249683673381
$ws.Cells.Item(15, 4).Value = 7.874172184888872
$ws.Cells.Item(15, 5).Value = 12.92760305708897
$ws.Cells.Item(15, 6).Value = 38.44821351043175
$ws.Cells.Item(15, 9).Value = 28.97012434695554
$ws.Cells.Item(15, 10).Value = 10.29230776523476
$ws.Cells.Item(15, 11).Value = 10.74775153795853
$ws.Cells.Item(15, 12).Value = 11.12169310094379
$ws.Cells.Item(15, 13).Value = 15.14319856695765
$ws.Cells.Item(15, 15).Value = 29.71454279507868
$ws.Cells.Item(16, 2).Value = 13.36058569187149
$ws.Cells.Item(16, 3).Value = 7.442738459637624
$ws.Cells.Item(16, 4).Value = 7.86122939738724
$ws.Cells.Item(16, 5).Value = 12.93984147800468
$ws.Cells.Item(16, 6).Value = 38.49746619423983
$ws.Cells.Item(16, 9).Value = 29.02848080036546
$ws.Cells.Item(16, 10).Value = 10.30301622908567
$ws.Cells.Item(16, 11).Value = 10.63649867866351
$ws.Cells.Item(16, 12).Value = 11.12254099853925
$ws.Cells.Item(16, 13).Value = 15.11445096060717
$ws.Cells.Item(16, 15).Value = 29.76572943584431
$ws.Cells.Item(17, 2).Value = 13.26712910506123
$ws.Cells.Item(17, 3).Value = 7.419523057965868
$ws.Cells.Item(17, 4).Value = 7.853519424105674
$ws.Cells.Item(17, 5).Value = 12.94767736061887
$ws.Cells.Item(17, 6).Value = 38.52956126717091
$ws.Cells.Item(17, 9).Value = 29.06545111709655
$ws.Cells.Item(17, 10).Value = 10.30973309162708
$ws.Cells.Item(17, 11).Value = 10.56797453476496
$ws.Cells.Item(17, 12).Value = 11.12336028020539
$ws.Cells.Item(17, 13).Value = 15.09722849342091
$ws.Cells.Item(17, 15).Value = 29.79845552319811
$ws.Cells.Item(18, 2).Value = 13.21323398082176
$ws.Cells.Item(18, 3).Value = 7.406096996168545
$ws.Cells.Item(18, 4).Value = 7.849170244135252
$ws.Cells.Item(18, 5).Value = 12.95230506323733
$ws.Cells.Item(18, 6).Value = 38.54871213701271
$ws.Cells.Item(18, 9).Value = 29.08714569049139
$ws.Cells.Item(18, 10).Value = 10.31365080029094
$ws.Cells.Item(18, 11).Value = 10.52847315778884
$ws.Cells.Item(18, 12).Value = 11.12394172621634
$ws.Cells.Item(18, 13).Value = 15.08747633252382
$ws.Cells.Item(18, 15).Value = 29.81776536157011
$ws.Cells.Item(19, 2).Value = 13.19496347420895
$ws.Cells.Item(19, 3).Value = 7.401538746899264
$ws.Cells.Item(19, 4).Value = 7.847712433613379
$ws.Cells.Item(19, 5).Value = 12.95389267183196
$ws.Cells.Item(19, 6).Value = 38.55531487752615
$ws.Cells.Item(19, 9).Value = 29.09456499797229
$ws.Cells.Item(19, 10).Value = 10.31498661437899
$ws.Cells.Item(19, 11).Value = 10.51508475272719
$ws.Cells.Item(19, 12).Value = 11.12415755070311
$ws.Cells.Item(19, 13).Value = 15.08420099993937
$ws.Cells.Item(19, 15).Value = 29.82438690785374
$ws.Cells.Item(20, 2).Value = 13.27709276041945
$ws.Cells.Item(20, 3).Value = 7.422001986891489
$ws.Cells.Item(20, 4).Value = 7.854331345395933
$ws.Cells.Item(20, 5).Value = 12.94683072827719
$ws.Cells.Item(20, 6).Value = 38.52607320651589
$ws.Cells.Item(20, 9).Value = 29.06147103851814
$ws.Cells.Item(20, 10).Value = 10.30901244807087
$ws.Cells.Item(20, 11).Value = 10.57527845270853
$ws.Cells.Item(20, 12).Value = 11.12326166569457
$ws.Cells.Item(20, 13).Value = 15.09904598410349
$ws.Cells.Item(20, 15).Value = 29.7949214038997
$ws.Cells.Item(21, 2).Value = 13.5512206604633
$ws.Cells.Item(21, 3).Value = 7.489870138270736
$ws.Cells.Item(21, 4).Value = 7.877584150175957
$ws.Cells.Item(21, 5).Value = 12.92454460656146
$ws.Cells.Item(21, 6).Value = 38.43607894126293
$ws.Cells.Item(21, 9).Value = 28.95541968878573
$ws.Cells.Item(21, 10).Value = 10.28958857738163
$ws.Cells.Item(21, 11).Value = 10.77637993736919
$ws.Cells.Item(21, 12).Value = 11.12156642713083
$ws.Cells.Item(21, 13).Value = 15.15074684538825
$ws.Cells.Item(21, 15).Value = 29.70173711188956
$ws.Cells.Item(22, 2).Value = 13.72874294267789
$ws.Cells.Item(22, 3).Value = 7.533540431564459
$ws.Cells.Item(22, 4).Value = 7.893522957543499
$ws.Cells.Item(22, 5).Value = 12.9110506581749
$ws.Cells.Item(22, 6).Value = 38.3834187595388
$ws.Cells.Item(22, 9).Value = 28.88993884420178
$ws.Cells.Item(22, 10).Value = 10.27737543802378
$ws.Cells.Item(22, 11).Value = 10.90676190826765
$ws.Cells.Item(22, 12).Value = 11.12143370158036
$ws.Cells.Item(22, 13).Value = 15.18586508164671
$ws.Cells.Item(22, 15).Value = 29.64517318261362
$ws.Cells.Item(23, 2).Value = 13.63416344596824
$ws.Cells.Item(23, 3).Value = 7.510297116177074
$ws.Cells.Item(23, 4).Value = 7.884949170079168
$ws.Cells.Item(23, 5).Value = 12.91815478120796
$ws.Cells.Item(23, 6).Value = 38.41096092143824
$ws.Cells.Item(23, 9).Value = 28.92453682900491
$ws.Cells.Item(23, 10).Value = 10.28384988354909
$ws.Cells.Item(23, 11).Value = 10.83728315430318
$ws.Cells.Item(23, 12).Value = 11.12141561899271
$ws.Cells.Item(23, 13).Value = 15.167002233939
$ws.Cells.Item(23, 15).Value = 29.67496536745254
$ws.Cells.Item(24, 2).Value = 13.27258870364895
$ws.Cells.Item(24, 3).Value = 7.420881509507038
$ws.Cells.Item(24, 4).Value = 7.853964015849338
$ws.Cells.Item(24, 5).Value = 12.94721310837437
$ws.Cells.Item(24, 6).Value = 38.52764798152172
$ws.Cells.Item(24, 9).Value = 29.06326906192988
$ws.Cells.Item(24, 10).Value = 10.30933807626439
$ws.Cells.Item(24, 11).Value = 10.57197667837718
$ws.Cells.Item(24, 12).Value = 11.12330590527223
$ws.Cells.Item(24, 13).Value = 15.09822383118113
$ws.Cells.Item(24, 15).Value = 29.79651763683093
$ws.Cells.Item(25, 2).Value = 12.87741170337103
$ws.Cells.Item(25, 3).Value = 7.321655085584034
$ws.Cells.Item(25, 4).Value = 7.823816014186607
$ws.Cells.Item(25, 5).Value = 12.98317432365396
$ws.Cells.Item(25, 6).Value = 38.67991272774468
$ws.Cells.Item(25, 9).Value = 29.22936381494213
$ws.Cells.Item(25, 10).Value = 10.33891057543492
$ws.Cells.Item(25, 11).Value = 10.2826017576847
$ws.Cells.Item(25, 12).Value = 11.12956192883052
$ws.Cells.Item(25, 13).Value = 15.02993462842789
$ws.Cells.Item(25, 15).Value = 29.94623505826386
